$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Mapping of row -> new "PERIOD TO EXPIRE" (column H) value.
# Each value is the previous value minus 1 (one day has passed).
$updates = @{
    3  = 377
    4  = 164
    5  = 190
    6  = 191
    7  = 359
    8  = 476
    9  = 170
    10 = 198
    11 = 189
    12 = 661
    13 = 161
    14 = 300
    15 = 657
    16 = 310
    17 = 346
    18 = 320
    19 = 608
    20 = 377
}

foreach ($row in $updates.Keys) {
    # Update PERIOD TO EXPIRE (column H)
    $ws.Cells.Item($row, 8).Value = $updates[$row]

    # Update LAST UPDATE (column I) to the new date, keeping it as plain text
    # (forcing a Text number format before assignment prevents Excel's COM
    # layer from auto-converting the date-like string into a real date value).
    $cellI = $ws.Cells.Item($row, 9)
    $cellI.NumberFormat = "@"
    $cellI.Value = "04-Nov-2025"
}
